$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 2
$ws.Range("H2").Value = 446
$ws.Range("I2").Value = 496
$ws.Range("J2").Value = 396
$ws.Range("K2").Value = 496
$ws.Range("L2").Value = 396
$ws.Range("M2").Value = -383
$ws.Range("N2").Value = -622

# Row 38
$ws.Range("H38").Value = 316.1111
$ws.Range("I38").Value = 123.28571
$ws.Range("J38").Value = 991
$ws.Range("K38").Value = 369.85713
$ws.Range("L38").Value = 2973
$ws.Range("M38").Value = 2.142870000000016
$ws.Range("N38").Value = -3717

# Row 40
$ws.Range("H40").Value = 2083.5
$ws.Range("I40").Value = 1440
$ws.Range("K40").Value = 1440
$ws.Range("M40").Value = -1265

# Row 41
$ws.Range("H41").Value = 62790.75
$ws.Range("I41").Value = 137.75
$ws.Range("J41").Value = 83675.086
$ws.Range("K41").Value = 137.75
$ws.Range("L41").Value = 83675.086
$ws.Range("M41").Value = 302.25
$ws.Range("N41").Value = -84555.086

# Row 125
$ws.Range("H125").Value = 1358.7
$ws.Range("I125").Value = 977.25
$ws.Range("J125").Value = 1613
$ws.Range("K125").Value = 8795.25
$ws.Range("L125").Value = 14517
$ws.Range("M125").Value = -6335.25
$ws.Range("N125").Value = -19437

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 2602.75
$ws.Range("I2").Value = 2470.3333
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 2470.3333
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -2357.3333
$ws.Range("N2").Value = -3226

# Row 32
$ws.Range("H32").Value = 4348.3667
$ws.Range("I32").Value = 3319.0476
$ws.Range("K32").Value = 3319.0476
$ws.Range("M32").Value = -3032.0476

# Row 44
$ws.Range("H44").Value = 75000
$ws.Range("J44").Value = 75000
$ws.Range("L44").Value = 75000
$ws.Range("N44").Value = -75976

# Row 55
$ws.Range("H55").Value = 79999.5
$ws.Range("I55").Value = 59999
$ws.Range("K55").Value = 59999
$ws.Range("M55").Value = -59684

# Row 61
$ws.Range("H61").Value = 20835046
$ws.Range("I61").Value = 30303836
$ws.Range("J61").Value = 3705
$ws.Range("K61").Value = 30303836
$ws.Range("L61").Value = 3705
$ws.Range("M61").Value = -30303624
$ws.Range("N61").Value = -4129

# Row 62
$ws.Range("H62").Value = 45999
$ws.Range("J62").Value = 45999
$ws.Range("L62").Value = 45999
$ws.Range("N62").Value = -47247

# Row 65
$ws.Range("H65").Value = 45999
$ws.Range("J65").Value = 45999
$ws.Range("L65").Value = 137997
$ws.Range("N65").Value = -144237

# Row 74
$ws.Range("H74").Value = 31254154
$ws.Range("I74").Value = 41670292
$ws.Range("K74").Value = 41670292
$ws.Range("M74").Value = -41669418

# Row 77
$ws.Range("H77").Value = 31254154
$ws.Range("I77").Value = 41670292
$ws.Range("K77").Value = 208351460
$ws.Range("M77").Value = -208347092

# Row 88
$ws.Range("H88").Value = 11906977
$ws.Range("J88").Value = 2779
$ws.Range("L88").Value = 2779
$ws.Range("N88").Value = -3591

# Row 91
$ws.Range("H91").Value = 11906977
$ws.Range("J91").Value = 2779
$ws.Range("L91").Value = 2779
$ws.Range("N91").Value = -5587

# Row 110
$ws.Range("H110").Value = 333334660
$ws.Range("I110").Value = 333334660
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 333334660
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -333332615
$ws.Range("N110").ClearContents()

# Row 116
$ws.Range("H116").Value = 2602.75
$ws.Range("I116").Value = 2470.3333
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2470.3333
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -176.3332999999998
$ws.Range("N116").Value = -7588

# Row 136
$ws.Range("H136").Value = 20835046
$ws.Range("I136").Value = 30303836
$ws.Range("J136").Value = 3705
$ws.Range("K136").Value = 90911508
$ws.Range("L136").Value = 11115
$ws.Range("M136").Value = -90908958
$ws.Range("N136").Value = -16215

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 2602.75
$ws.Range("I3").Value = 2470.3333
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 2470.3333
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -2356.3333
$ws.Range("N3").Value = -3228

# Row 137
$ws.Range("H137").Value = 54998.5
$ws.Range("J137").Value = 54998.5
$ws.Range("L137").Value = 54998.5
$ws.Range("N137").Value = -65198.5

# Row 138
$ws.Range("H138").Value = 53566
$ws.Range("J138").Value = 59994.5
$ws.Range("L138").Value = 59994.5
$ws.Range("N138").Value = -70274.5

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 546
$ws.Range("J22").Value = 139.75
$ws.Range("L22").Value = 139.75
$ws.Range("N22").Value = -839.75

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 9796090
$ws.Range("I4").Value = 10728642
$ws.Range("J4").Value = 4284.8335
$ws.Range("K4").Value = 32185926
$ws.Range("L4").Value = 12854.5005
$ws.Range("M4").Value = -32185814
$ws.Range("N4").Value = -13078.5005

# Row 10
$ws.Range("H10").Value = 627.05554
$ws.Range("I10").Value = 254.77777
$ws.Range("K10").Value = 764.33331
$ws.Range("M10").Value = -625.33331

# Row 38
$ws.Range("H38").Value = 94.08
$ws.Range("I38").Value = 100.64286
$ws.Range("J38").Value = 85.72727
$ws.Range("K38").Value = 301.92858
$ws.Range("L38").Value = 257.18181
$ws.Range("M38").Value = 45.07141999999999
$ws.Range("N38").Value = -951.18181

# Row 70
$ws.Range("H70").Value = 3353
$ws.Range("I70").Value = 3353
$ws.Range("K70").Value = 10059
$ws.Range("M70").Value = -9744

# Row 73
$ws.Range("H73").Value = 3353
$ws.Range("I73").Value = 3353
$ws.Range("K73").Value = 10059
$ws.Range("M73").Value = -8967

# Row 75
$ws.Range("H75").Value = 806
$ws.Range("J75").Value = 1141.3334
$ws.Range("L75").Value = 3424.0002
$ws.Range("N75").Value = -5420.0002

# Row 78
$ws.Range("H78").Value = 806
$ws.Range("J78").Value = 1141.3334
$ws.Range("L78").Value = 10272.0006
$ws.Range("N78").Value = -20256.0006

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 99
$ws.Range("H99").Value = 1452.6666
$ws.Range("I99").Value = 1452.6666
$ws.Range("K99").Value = 1452.6666
$ws.Range("M99").Value = 793.3334

# Row 102
$ws.Range("H102").Value = 2653.682
$ws.Range("J102").Value = 3520.2727
$ws.Range("L102").Value = 3520.2727
$ws.Range("N102").Value = -6764.2727

# Row 113
$ws.Range("H113").Value = 2052.2693
$ws.Range("I113").Value = 1573.6666
$ws.Range("K113").Value = 1573.6666
$ws.Range("M113").Value = 596.3334

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 2437
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value = 2437
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Row 55
$ws.Range("H55").Value = 692.3333
$ws.Range("J55").Value = 743.8570999999999
$ws.Range("L55").Value = 743.8570999999999
$ws.Range("N55").Value = -1089.8571

# Row 61
$ws.Range("H61").Value = 1725
$ws.Range("I61").Value = 1450
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1450
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1248
$ws.Range("N61").Value = -2404

# Row 113
$ws.Range("H113").Value = 1725
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1450
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 720
$ws.Range("N113").Value = -6340

# Row 136
$ws.Range("H136").Value = 4332.8335
$ws.Range("I136").Value = 3999
$ws.Range("J136").Value = 4499.75
$ws.Range("K136").Value = 11997
$ws.Range("L136").Value = 13499.25
$ws.Range("M136").Value = -9447
$ws.Range("N136").Value = -18599.25

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 130
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040

# Row 132
$ws.Range("H132").Value = 3184.5
$ws.Range("I132").Value = 2627.95
$ws.Range("K132").Value = 7883.849999999999
$ws.Range("M132").Value = -5353.849999999999

# Row 136
$ws.Range("H136").Value = 4817.143
$ws.Range("I136").Value = 1888.8
$ws.Range("K136").Value = 5666.4
$ws.Range("M136").Value = -3116.4
